$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after row 19 so the existing row 19 data moves to row 20.
$ws.Rows.Item(20).Insert()

# Copy the (old) row 19 values down into the newly created row 20.
for ($col = 1; $col -le 20; $col++) {
    $src = $ws.Cells.Item(19, $col)
    $dst = $ws.Cells.Item(20, $col)
    $dst.Value = $src.Value2
}

# D20 keeps the "date" number format, matching D19's original style.
$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat

# Update row 19 with the new weekly values.
$ws.Range("D19").Value = 44448
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 22000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 22000
$ws.Range("S19").Value = 1100
